# ---------------------------------------------------------------------------
# feat: update package for D2.35
#
# Updates the "Package info" metadata (version/DHIS2 version/created/
# identifier), re-labels a few Tracked Entity Attribute names in the
# programTrackedEntityAttributes sheet, re-orders the dataElementGroups
# listing, reorders/extends the trackedEntityAttributes sheet (alphabetical,
# with three newly-documented attributes), and bumps the "Last updated" date
# on the programs sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Package info" sheet
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")

# Column B is narrower in the new version of the package.
$wsInfo.Columns.Item(2).ColumnWidth = 56.83333333333333

$wsInfo.Range("B4").Value = "V1.1.2"
$wsInfo.Range("B5").Value = "DHIS2.35.3-3492688"
$wsInfo.Range("B6").Value = "20210408T081801"
$wsInfo.Range("B7").Value = "COVAC_TRACKER_V1.1.2_DHIS2.35.3-3492688_20210408T081801"

# ---------------------------------------------------------------------------
# 2. "programTrackedEntityAttributes" sheet - a handful of rows had the raw
#    UID in the "Tracked Entity Attribute Name" column; they now show the
#    human readable attribute name instead.
# ---------------------------------------------------------------------------
$wsPTEA = $wb.Worksheets.Item("programTrackedEntityAttributes")

$wsPTEA.Range("B4").Value = "First Name"
$wsPTEA.Range("B5").Value = "Surname"
$wsPTEA.Range("B6").Value = "Sex"
$wsPTEA.Range("B8").Value = "Date of birth"
$wsPTEA.Range("B10").Value = "Home Address"

# ---------------------------------------------------------------------------
# 3. "dataElementGroups" sheet - the "Data Element" column was re-ordered.
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")

$dataElementOrder = @(
    "COVAC - Dose Expiry Date",
    "COVAC - Underlying condition Other",
    "COVAC - Dose Number",
    "COVAC- Batch Number",
    "COVAC - Renal Disease",
    "COVAC - Malignancy",
    "COVAC - Vaccine Name",
    "COVAC - Pregnancy",
    "COVAC - Immunodeficiency",
    "COVAC Previously infected with COVID",
    "COVAC Suggested date for next dose",
    "COVAC - Multiple products used - Explain",
    "COVAC - AEFIs present",
    "COVAC - Cardiovascular Disease",
    "COVAC - Vaccine Manufacturer",
    "COVAC - Last Dose",
    "COVAC - Pregnancy gestation",
    "COVAC - Allergic reaction after first dose",
    "COVAC - Chronic Lung Disease",
    "COVAC - Diabetes",
    "COVAC - Neurological/Neuromuscular",
    "COVAC - Underlying condition",
    "COVAC - Total doses"
)

for ($i = 0; $i -lt $dataElementOrder.Length; $i++) {
    $rowNum = $i + 2
    $wsDEG.Cells.Item($rowNum, 2).Value = $dataElementOrder[$i]
}

# ---------------------------------------------------------------------------
# 4. "trackedEntityAttributes" sheet - the rows are now listed alphabetically
#    and three additional attributes (First Name, Home Address, Sex, Surname,
#    Date of birth) are fully documented. Column B is a little wider to fit
#    the new "Code" values.
# ---------------------------------------------------------------------------
$wsTEA = $wb.Worksheets.Item("trackedEntityAttributes")

$wsTEA.Columns.Item(2).ColumnWidth = 21.833333333333332

# Rows 1-3 (header, "Area Urban Rural", "COVID - Occupation") are unchanged.
# Rebuild rows 4-12 from scratch so no stale values survive from the old
# (shorter) layout.
$teaRows = @(
    @(4,  "4", "Date of birth",                   "patinfo_ageonsetunit", "",                                     "NI0QRzJvQ0k"),
    @(5,  "5", "Date of birth is estimated",       "",                     "",                                     "Z1rLc1rVHK8"),
    @(6,  "4", "First Name",                       "first_name",           "",                                     "sB1IHYu2xQT"),
    @(7,  "5", "Home Address",                     "patinfo_resadmin0",    "",                                     "Xhdn49gUd52"),
    @(8,  "4", "Mobile phone number",              "",                     "",                                     "fctSQp5nAYl"),
    @(9,  "5", "National ID",                      "",                     "",                                     "Ewi7FUfcHAD"),
    @(10, "4", "Sex",                               "patinfo_sex",          "",                                     "oindugucx72"),
    @(11, "5", "Surname",                           "surname",              "The patient's surname (family name)", "ENRjVGxVL6l"),
    @(12, "4", "Unique System Identifier (EPI)",    "",                     "System-generated unique ID following pattern: EPI prefix + value randomly generated (#####) - Customize the length depending on the target population of your implementation", "KSr2yTdu1AI")
)

foreach ($row in $teaRows) {
    $rowNum = $row[0]
    $styleSourceRow = if ($row[1] -eq "4") { 2 } else { 3 }

    # Copy the banding/fill format from an existing row with the matching
    # style (even rows use the "4" style, odd rows use the "5" style) onto
    # this row, then fill in the values.
    $wsTEA.Range("A" + $styleSourceRow + ":E" + $styleSourceRow).Copy() | Out-Null
    $wsTEA.Range("A" + $rowNum + ":E" + $rowNum).PasteSpecial(-4122) | Out-Null

    $wsTEA.Cells.Item($rowNum, 1).Value = $row[2]
    $wsTEA.Cells.Item($rowNum, 2).Value = $row[3]
    $wsTEA.Cells.Item($rowNum, 3).Value = $row[4]
    $wsTEA.Cells.Item($rowNum, 4).Value = ""
    $wsTEA.Cells.Item($rowNum, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# 5. "programs" sheet - "Last updated" date bump.
# ---------------------------------------------------------------------------
$wsPrograms = $wb.Worksheets.Item("programs")
$wsPrograms.Range("C2").Value = "2021-03-19"
